$wb = $excel.ActiveWorkbook

# Fix the typo in Sheet2!A2: "countrys" -> "country"
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A2").Value = "country"
$ws2.Range("B2").Select()

# Change Sheet1!B3 from "basketball" to "golf", then make Sheet1 the active
# sheet with B3 selected (this also activates the sheet / sets tabSelected).
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B3").Value = "golf"
$ws1.Range("B3").Select()
